$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Reboot")
$ws2 = $wb.Worksheets.Item("City")

# --- Update existing City values (keep them as text, matching original shared-string type) ---
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "1000"
$ws2.Range("B2").NumberFormat = "General"

$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "5000"
$ws2.Range("C2").NumberFormat = "General"

# --- Add the two new columns (headers) on the City sheet ---
$ws2.Range("J1").Value = "UpdatedByEmailId"
$ws2.Range("K1").Value = "TONS/YEAR/CAPITA"

# --- New data row values ---
# J2: email address, rendered as a mailto hyperlink using the same visual style as
# the existing hyperlink cell on the Reboot sheet (blue text, no underline).
$ws2.Range("J2").Value = "narayanarc21@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("J2"), "mailto:narayanarc21@gmail.com", "", "", "narayanarc21@gmail.com")
$ws2.Range("J2").Font.Underline = -4142
$ws2.Range("J2").Font.Name = "Arial"
$ws2.Range("J2").Font.Color = 16711680

# K2: numeric-looking value stored as text, like the rest of the row.
$ws2.Range("K2").NumberFormat = "@"
$ws2.Range("K2").Value = "15"
$ws2.Range("K2").NumberFormat = "General"

# --- Column width for the newly added column J ---
$ws2.Columns.Item(10).ColumnWidth = 22

# --- Selection state: Reboot keeps a non-active selection on H2 ... ---
$ws1.Activate()
$ws1.Range("H2").Select()

# --- ... while City becomes the active / selected sheet with C2 selected ---
$ws2.Activate()
$ws2.Range("C2").Select()
